$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IFStFS")

# Copy formatting from the existing "hydrogen if" header/data column (K)
# into the two new columns (L, M) that will hold the new subscript
# elements: "green hydrogen if" and "low carbon hydrogen if".
$ws.Range("K1").Copy() | Out-Null
$ws.Range("L1:M1").PasteSpecial(-4122) | Out-Null

$ws.Range("K2:K26").Copy() | Out-Null
$ws.Range("L2:M26").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# New header labels
$ws.Range("L1").Value = "green hydrogen if"
$ws.Range("M1").Value = "low carbon hydrogen if"

# New columns start at zero for every industry, same as the other
# fuel-shifting flags feeding the demand structure.
$ws.Range("L2:M26").Value = 0

# Match column widths of neighboring fuel columns
$ws.Range("L:M").ColumnWidth = $ws.Range("K:K").ColumnWidth

# Make the IFStFS sheet the active tab with M2 selected
$ws.Activate()
$ws.Range("M2").Select() | Out-Null
